$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 72
$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 22.5
